# ExportToExcel: turn the generic "Data" tab into a real "Conversations"
# export sheet (headers, autofilter, styling) and add a blank "SURA" tab,
# plus point Main's summary row at the renamed sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "Data" -> "Conversations", then append a new blank "SURA"
#    sheet right after it.
# ---------------------------------------------------------------------
$wsConv = $wb.Worksheets.Item("Data")
$wsConv.Name = "Conversations"

$wsSura = $wb.Worksheets.Add($null, $wsConv)
$wsSura.Name = "SURA"
# Match the plain default page setup used elsewhere in the workbook.
$wsSura.PageSetup.LeftMargin = 54
$wsSura.PageSetup.RightMargin = 54
$wsSura.PageSetup.TopMargin = 72
$wsSura.PageSetup.BottomMargin = 72
$wsSura.PageSetup.HeaderMargin = 36
$wsSura.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Main sheet: refresh the summary blurb + label, and re-point the
#    COUNTA formula (now in B2) at the renamed Conversations sheet.
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Range("A1").Value = "This is a test - summary stuff here…"
$wsMain.Range("A1:J1").HorizontalAlignment = -4131   # xlLeft

$wsMain.Range("A2").Value = "Message Lines"
$wsMain.Range("B2").Formula = "=COUNTA(Conversations!A:A)"

# ---------------------------------------------------------------------
# 3. Conversations sheet: write the export header row, style it, size
#    column A, and turn on an AutoFilter over the header.
# ---------------------------------------------------------------------
$headers = @("assessmentId", "createdBy", "conversationId", "eventKey", "text", "mcsRawScore", "agentId", "time", "note")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsConv.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $wsConv.Range("A1:I1")
$headerRange.Font.Color = 0   # solid black text for the header labels

$wsConv.Columns.Item(1).ColumnWidth = 12.17   # renders as width 13 in xml

$wsConv.Range("A1:I1").AutoFilter() | Out-Null
$filterName = $wsConv.Names.Add("_xlnm._FilterDatabase", "=Conversations!`$A`$1:`$I`$1")
$filterName.Visible = $false

# Leave Conversations as the active/selected sheet with the header row
# highlighted, mirroring the authored workbook view.
$wsConv.Range("A1:I1").Select() | Out-Null
$wsConv.Activate()
